$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update TPM-derived numeric values (columns E:T) for all 20 data rows (rows 2-21).
# Columns A-D (Sending cluster, Ligand, Receptor, Target cluster) are unchanged.
$updates = @(
    @(2, 7, 9.912502333333334),
    @(2, 8, 29.737507),
    @(2, 9, 0.306800202315277),
    @(2, 10, 0.3105483022825659),
    @(2, 13, 6.712486666666666),
    @(2, 14, 20.13746),
    @(2, 15, 0.6330487633990675),
    @(2, 16, 0.6414503882251803),
    @(2, 17, 66.53753974580222),
    @(2, 18, 598.8378577122199),
    @(2, 19, 0.1942194886862698),
    @(2, 20, 0.1992013290618225),
    @(3, 7, 9.912502333333334),
    @(3, 8, 29.737507),
    @(3, 9, 0.306800202315277),
    @(3, 10, 0.3105483022825659),
    @(3, 15, 0.290741083484562),
    @(3, 16, 0.2945997080427384),
    @(3, 17, 30.55877764332134),
    @(3, 18, 275.028998789892),
    @(3, 19, 0.08919942323442646),
    @(3, 20, 0.09148743918561197),
    @(4, 7, 9.912502333333334),
    @(4, 8, 29.737507),
    @(4, 9, 0.306800202315277),
    @(4, 10, 0.3105483022825659),
    @(4, 13, 0.2495096666666667),
    @(4, 14, 0.748529),
    @(4, 15, 0.02353103905946135),
    @(4, 16, 0.02384333563656022),
    @(4, 17, 2.473265153022556),
    @(4, 18, 22.259386377203),
    @(4, 19, 0.007219327544131429),
    @(4, 20, 0.00740450740268718),
    @(5, 7, 9.912502333333334),
    @(5, 8, 29.737507),
    @(5, 9, 0.306800202315277),
    @(5, 10, 0.3105483022825659),
    @(5, 13, 0.4166465),
    @(5, 14, 0.8332930000000001),
    @(5, 15, 0.03929356804674715),
    @(5, 16, 0.02654337331298611),
    @(5, 17, 4.130009403425167),
    @(5, 18, 24.780056420551),
    @(5, 19, 0.01205527462643113),
    @(5, 20, 0.008242999519200203),
    @(6, 7, 9.912502333333334),
    @(6, 8, 29.737507),
    @(6, 9, 0.306800202315277),
    @(6, 10, 0.3105483022825659),
    @(6, 11, 2.0),
    @(6, 12, 0.6666666666666666),
    @(6, 13, 0.1419326666666667),
    @(6, 14, 0.425798),
    @(6, 15, 0.01338554601016197),
    @(6, 16, 0.01356319478253491),
    @(6, 17, 1.406907889509556),
    @(6, 18, 12.662171005586),
    @(6, 19, 0.004106688224018141),
    @(6, 20, 0.00421202711324397),
    @(7, 9, 0.648195099606228),
    @(7, 10, 0.6561139341222959),
    @(7, 13, 6.712486666666666),
    @(7, 14, 20.13746),
    @(7, 15, 0.6330487633990675),
    @(7, 16, 0.6414503882251803),
    @(7, 17, 140.5778316885289),
    @(7, 18, 1265.20048519676),
    @(7, 19, 0.410339106247058),
    @(7, 20, 0.4208645377626971),
    @(8, 9, 0.648195099606228),
    @(8, 10, 0.6561139341222959),
    @(8, 15, 0.290741083484562),
    @(8, 16, 0.2945997080427384),
    @(8, 19, 0.1884569455688983),
    @(8, 20, 0.1932909734352009),
    @(9, 9, 0.648195099606228),
    @(9, 10, 0.6561139341222959),
    @(9, 13, 0.2495096666666667),
    @(9, 14, 0.748529),
    @(9, 15, 0.02353103905946135),
    @(9, 16, 0.02384333563656022),
    @(9, 17, 5.225414912108222),
    @(9, 18, 47.02873420897399),
    @(9, 19, 0.01525270420698559),
    @(9, 20, 0.01564394474710187),
    @(10, 9, 0.648195099606228),
    @(10, 10, 0.6561139341222959),
    @(10, 13, 0.4166465),
    @(10, 14, 0.8332930000000001),
    @(10, 15, 0.03929356804674715),
    @(10, 16, 0.02654337331298611),
    @(10, 17, 8.725717377059667),
    @(10, 18, 52.354304262358),
    @(10, 19, 0.02546989825394537),
    @(10, 20, 0.01741547708926008),
    @(11, 9, 0.648195099606228),
    @(11, 10, 0.6561139341222959),
    @(11, 11, 2.0),
    @(11, 12, 0.6666666666666666),
    @(11, 13, 0.1419326666666667),
    @(11, 14, 0.425798),
    @(11, 15, 0.01338554601016197),
    @(11, 16, 0.01356319478253491),
    @(11, 17, 2.972458273154222),
    @(11, 18, 26.752124458388),
    @(11, 19, 0.008676445329340685),
    @(11, 20, 0.008899001088035974),
    @(12, 5, 2.0),
    @(12, 6, 1.0),
    @(12, 7, 1.169852),
    @(12, 8, 2.339704),
    @(12, 9, 0.03620789364881174),
    @(12, 10, 0.0244334908452053),
    @(12, 13, 6.712486666666666),
    @(12, 14, 20.13746),
    @(12, 15, 0.6330487633990675),
    @(12, 16, 0.6414503882251803),
    @(12, 17, 7.852615951973333),
    @(12, 18, 47.11569571184),
    @(12, 19, 0.02292136229966522),
    @(12, 20, 0.01567287218835333),
    @(13, 5, 2.0),
    @(13, 6, 1.0),
    @(13, 7, 1.169852),
    @(13, 8, 2.339704),
    @(13, 9, 0.03620789364881174),
    @(13, 10, 0.0244334908452053),
    @(13, 15, 0.290741083484562),
    @(13, 16, 0.2945997080427384),
    @(13, 17, 3.606480577904001),
    @(13, 18, 21.63888346742401),
    @(13, 19, 0.01052712223014931),
    @(13, 20, 0.007198099269462402),
    @(14, 5, 2.0),
    @(14, 6, 1.0),
    @(14, 7, 1.169852),
    @(14, 8, 2.339704),
    @(14, 9, 0.03620789364881174),
    @(14, 10, 0.0244334908452053),
    @(14, 13, 0.2495096666666667),
    @(14, 14, 0.748529),
    @(14, 15, 0.02353103905946135),
    @(14, 16, 0.02384333563656022),
    @(14, 17, 0.2918893825693333),
    @(14, 18, 1.751336295416),
    @(14, 19, 0.0008520093597110117),
    @(14, 20, 0.0005825759229950515),
    @(15, 5, 2.0),
    @(15, 6, 1.0),
    @(15, 7, 1.169852),
    @(15, 8, 2.339704),
    @(15, 9, 0.03620789364881174),
    @(15, 10, 0.0244334908452053),
    @(15, 13, 0.4166465),
    @(15, 14, 0.8332930000000001),
    @(15, 15, 0.03929356804674715),
    @(15, 16, 0.02654337331298611),
    @(15, 17, 0.4874147413180001),
    @(15, 18, 1.949658965272),
    @(15, 19, 0.001422737332918968),
    @(15, 20, 0.0006485472688437128),
    @(16, 5, 2.0),
    @(16, 6, 1.0),
    @(16, 7, 1.169852),
    @(16, 8, 2.339704),
    @(16, 9, 0.03620789364881174),
    @(16, 10, 0.0244334908452053),
    @(16, 11, 2.0),
    @(16, 12, 0.6666666666666666),
    @(16, 13, 0.1419326666666667),
    @(16, 14, 0.425798),
    @(16, 15, 0.01338554601016197),
    @(16, 16, 0.01356319478253491),
    @(16, 17, 0.1660402139653334),
    @(16, 18, 0.9962412837920002),
    @(16, 19, 0.0004846624263672208),
    @(16, 20, 0.0003313961955508028),
    @(17, 5, 1.0),
    @(17, 6, 0.3333333333333333),
    @(17, 7, 0.2842186666666667),
    @(17, 8, 0.852656),
    @(17, 9, 0.00879680442968319),
    @(17, 10, 0.008904272749933054),
    @(17, 13, 6.712486666666666),
    @(17, 14, 20.13746),
    @(17, 15, 0.6330487633990675),
    @(17, 16, 0.6414503882251803),
    @(17, 17, 1.907814010417778),
    @(17, 18, 17.17032609376),
    @(17, 19, 0.005568806166074382),
    @(17, 20, 0.005711649212307451),
    @(18, 5, 1.0),
    @(18, 6, 0.3333333333333333),
    @(18, 7, 0.2842186666666667),
    @(18, 8, 0.852656),
    @(18, 9, 0.00879680442968319),
    @(18, 10, 0.008904272749933054),
    @(18, 15, 0.290741083484562),
    @(18, 16, 0.2945997080427384),
    @(18, 17, 0.8762040849706668),
    @(18, 18, 7.885836764736),
    @(18, 19, 0.002557592451087885),
    @(18, 20, 0.002623196152463189),
    @(19, 5, 1.0),
    @(19, 6, 0.3333333333333333),
    @(19, 7, 0.2842186666666667),
    @(19, 8, 0.852656),
    @(19, 9, 0.00879680442968319),
    @(19, 10, 0.008904272749933054),
    @(19, 13, 0.2495096666666667),
    @(19, 14, 0.748529),
    @(19, 15, 0.02353103905946135),
    @(19, 16, 0.02384333563656022),
    @(19, 17, 0.07091530478044444),
    @(19, 18, 0.6382377430239999),
    @(19, 19, 0.0002069979486333178),
    @(19, 20, 0.0002123075637761309),
    @(20, 5, 1.0),
    @(20, 6, 0.3333333333333333),
    @(20, 7, 0.2842186666666667),
    @(20, 8, 0.852656),
    @(20, 9, 0.00879680442968319),
    @(20, 10, 0.008904272749933054),
    @(20, 13, 0.4166465),
    @(20, 14, 0.8332930000000001),
    @(20, 15, 0.03929356804674715),
    @(20, 16, 0.02654337331298611),
    @(20, 17, 0.1184187127013333),
    @(20, 18, 0.7105122762080001),
    @(20, 19, 0.0003456578334516832),
    @(20, 20, 0.0002363494356821225),
    @(21, 5, 1.0),
    @(21, 6, 0.3333333333333333),
    @(21, 7, 0.2842186666666667),
    @(21, 8, 0.852656),
    @(21, 9, 0.00879680442968319),
    @(21, 10, 0.008904272749933054),
    @(21, 11, 2.0),
    @(21, 12, 0.6666666666666666),
    @(21, 13, 0.1419326666666667),
    @(21, 14, 0.425798),
    @(21, 15, 0.01338554601016197),
    @(21, 16, 0.01356319478253491),
    @(21, 17, 0.04033991327644445),
    @(21, 18, 0.363059219488),
    @(21, 19, 0.0001177500304359209),
    @(21, 20, 0.0001207703857041597)
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

Write-Output ("Updated " + $updates.Count + " cells")